# Automatische test-sync: 2025-08-05 19:47:50
$wb = $excel.ActiveWorkbook

# --- 1) "Logs" sheet: append new row 53 with the new test mail entry ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A53").Value = "Wanneer wordt het geld van mijn retour overgemaakt?"
$logs.Range("B53").Value = "mailmind.test@zohomail.eu"
$logs.Range("C53").Value = "Testmail #12: Wanneer wordt het geld van mijn retour overgemaakt?"
$logs.Range("D53").Value = "Retour / Terugbetaling"
$logs.Range("E53").Value = "Geachte klant,
Dank u voor uw e-mail. Om uw terugbetaling te kunnen verwerken en het geld over te maken, hebben wij wat meer informatie nodig. Kunt u ons alstublieft uw ordernummer en de naam waaronder de bestelling geplaatst is, doorgeven? Zodra we deze gegevens ontvangen hebben, zullen we het proces voor de terugbetaling in gang zetten.
Met vriendelijke groet,
[Naam bedrijf] E-mailassistent"
$logs.Range("F53").Value = "2025-08-05 19:47:12"
$logs.Range("G53").Value = "Ja"
$logs.Range("H53").Value = "Nee"
$logs.Range("I53").Value = "Ja"
$logs.Range("J53").Value = "Nee"

# Writing the multi-line E53 text triggers an automatic row-height change;
# AutoFit() restores it back to the sheet's default height with no explicit
# custom height, just like all the other rows in this log.
$logs.Rows.Item(53).AutoFit()

# Extend the conditional formatting ranges (D,G,H,I,J) from row 52 to row 53
# so newly added row 53 keeps being highlighted like the rest of the log.
$cols = @("D","G","H","I","J")
foreach ($col in $cols) {
    $oldRange = $logs.Range($col + "2:" + $col + "52")
    $newRange = $logs.Range($col + "2:" + $col + "53")
    $fcs = $oldRange.FormatConditions
    $cnt = $fcs.Count()
    for ($i = 1; $i -le $cnt; $i++) {
        $fc = $fcs.Item($i)
        $fc.ModifyAppliesToRange($newRange)
    }
}

# --- 2) "Dashboard" sheet: swap the order/counts of the two categories ---
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A5").Value = "Retour / Terugbetaling"
$dash.Range("B5").Value = 5
$dash.Range("A6").Value = "Klantenservice / Contact"
$dash.Range("B6").Value = 4
